# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.941.14"
$ws.Range("E2").Value = "  +4.24%  "
$ws.Range("D3").Value = "'2.285.28"
$ws.Range("E3").Value = "  +5.07%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'252.28"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("E6").Value = "  +4.55%  "
$ws.Range("E7").Value = "  +10.08%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.658"
$ws.Range("E9").Value = "  +14.21%  "
$ws.Range("D10").Value = "'38.86"
$ws.Range("E10").Value = "  +7.26%  "
$ws.Range("D11").Value = "'0.0981"
$ws.Range("E11").Value = "  +5.00%  "
$ws.Range("D12").Value = "'59.95"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "'7.39"
$ws.Range("E13").Value = "  +8.11%  "
$ws.Range("D14").Value = "'0.105"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "'2.625.11"
$ws.Range("E15").Value = "  +5.09%  "
$ws.Range("D16").Value = "'15.06"
$ws.Range("E16").Value = "  +5.58%  "
$ws.Range("D17").Value = "'0.892"
$ws.Range("E17").Value = "  +5.50%  "
$ws.Range("D18").Value = "'2.276.91"
$ws.Range("E18").Value = "  +4.81%  "
$ws.Range("D19").Value = "'42.845.31"
$ws.Range("E19").Value = "  +4.28%  "
$ws.Range("E20").Value = "  +7.39%  "
$ws.Range("D21").Value = "'6.38"
$ws.Range("E21").Value = "  +5.42%  "
$ws.Range("D22").Value = "'73.52"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").Value = "'237.75"
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("E24").Value = "  +6.31%  "
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").Value = "'11.66"
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "'2.20"
$ws.Range("E30").Value = "  +4.09%  "
$ws.Range("D31").Value = "'167.92"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'21.11"
$ws.Range("E32").Value = "  +4.69%  "
$ws.Range("D33").Value = "'6.33"
$ws.Range("E33").Value = "  +11.36%  "
$ws.Range("E34").Value = "  +5.79%  "
$ws.Range("D35").Value = "'0.0814"
$ws.Range("E35").Value = "  +8.50%  "
$ws.Range("D36").Value = "'31.13"
$ws.Range("E36").Value = "  +27.93%  "
$ws.Range("E37").Value = "  +4.98%  "
$ws.Range("D38").Value = "'4.74"
$ws.Range("E38").Value = "  +20.74%  "
$ws.Range("D39").Value = "'4.80"
$ws.Range("E39").Value = "  +6.28%  "
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").Value = "'13.37"
$ws.Range("E41").Value = "  +17.87%  "
$ws.Range("D42").Value = "'2.33"
$ws.Range("E42").Value = "  +5.55%  "
$ws.Range("D43").Value = "'6.02"
$ws.Range("E43").Value = "  +10.03%  "
$ws.Range("D44").Value = "'0.216"
$ws.Range("E44").Value = "  +14.36%  "
$ws.Range("D45").Value = "'9.24"
$ws.Range("E45").Value = "  +9.18%  "
$ws.Range("D46").Value = "'5.00"
$ws.Range("E46").Value = "  -7.93%  "
$ws.Range("D47").Value = "'61.59"
$ws.Range("E48").Value = "  +4.70%  "
$ws.Range("E49").Value = "  +4.48%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  +5.67%  "
